$d = $word.ActiveDocument

# --- Image 1 (anchorId 1E20E740) ---
$shp1 = $d.InlineShapes.Item(1)
$rng1 = $shp1.Range
$xml1 = $rng1.WordOpenXML

$xml1 = $xml1.Replace(
    'wp14:anchorId="1E20E740" wp14:editId="10B4432C"><wp:extent cx="10488489" cy="4915586"/><wp:effectExtent l="0" t="0" r="8255" b="0"/>',
    'wp14:anchorId="1E20E740" wp14:editId="343E43E3"><wp:extent cx="5240434" cy="2456007"/><wp:effectExtent l="0" t="0" r="0" b="1905"/>'
)
$xml1 = $xml1.Replace(
    '<a:off x="0" y="0"/><a:ext cx="10488489" cy="4915586"/>',
    '<a:off x="0" y="0"/><a:ext cx="5264812" cy="2467432"/>'
)
$xml1 = $xml1.Replace(
    '<w:r w:rsidRPr="003E3668"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1E20E740"',
    '<w:r w:rsidRPr="003E3668"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1E20E740"'
)

$rng1.Text = ""
$null = $rng1.InsertXML($xml1)

# --- Image 2 (anchorId 04819F95) ---
$shp2 = $d.InlineShapes.Item(2)
$rng2 = $shp2.Range
$xml2 = $rng2.WordOpenXML

$xml2 = $xml2.Replace(
    'wp14:anchorId="04819F95" wp14:editId="066F753A"><wp:extent cx="10907647" cy="5515745"/><wp:effectExtent l="0" t="0" r="8255" b="8890"/>',
    'wp14:anchorId="04819F95" wp14:editId="4B8D2070"><wp:extent cx="5844226" cy="2955290"/><wp:effectExtent l="0" t="0" r="4445" b="0"/>'
)
$xml2 = $xml2.Replace(
    '<a:off x="0" y="0"/><a:ext cx="10907647" cy="5515745"/>',
    '<a:off x="0" y="0"/><a:ext cx="5864496" cy="2965540"/>'
)
$xml2 = $xml2.Replace(
    '<w:r w:rsidRPr="003E3668"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="04819F95"',
    '<w:r w:rsidRPr="003E3668"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="04819F95"'
)

$rng2.Text = ""
$null = $rng2.InsertXML($xml2)
